# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 57
    $ws.Range("F4").Value = 543
    $ws.Range("F9").Value = 4534
    $ws.Range("F10").Value = 4383
    $ws.Range("F12").Value = 18
    $ws.Range("F13").Value = 142
}
